# Update FFXIV Leve profit-tracking figures (currentAveragePrice / Leve
# price / profit columns H:N) across the eight job-table sheets, per the
# latest Universalis price pull from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 1510.7142
$ws.Range("J19").Value = 1444
$ws.Range("L19").Value = 1444
$ws.Range("N19").Value = -1794

# Row 39
$ws.Range("H39").Value = 2707.4
$ws.Range("J39").Value = 4549.8
$ws.Range("L39").Value = 13649.4
$ws.Range("N39").Value = -14241.4

# Row 70
$ws.Range("H70").Value = 3713.15
$ws.Range("I70").Value = 976.8
$ws.Range("K70").Value = 2930.4
$ws.Range("M70").Value = -2660.4

# Row 73
$ws.Range("H73").Value = 3713.15
$ws.Range("I73").Value = 976.8
$ws.Range("K73").Value = 2930.4
$ws.Range("M73").Value = -1994.4

# Row 88
$ws.Range("I88").Value = 66673360
$ws.Range("J88").Value = 2561000.8
$ws.Range("K88").Value = 66673360
$ws.Range("L88").Value = 2561000.8
$ws.Range("M88").Value = -66672954
$ws.Range("N88").Value = -2561812.8

# Row 91
$ws.Range("I91").Value = 66673360
$ws.Range("J91").Value = 2561000.8
$ws.Range("K91").Value = 66673360
$ws.Range("L91").Value = 2561000.8
$ws.Range("M91").Value = -66671956
$ws.Range("N91").Value = -2563808.8

# Row 116
$ws.Range("H116").Value = 6817
$ws.Range("I116").Value = 6777.067
$ws.Range("J116").Value = 7056.6
$ws.Range("K116").Value = 6777.067
$ws.Range("L116").Value = 7056.6
$ws.Range("M116").Value = -3335.067
$ws.Range("N116").Value = -13940.6

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 142858820
$ws.Range("I61").Value = 142858820
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 142858820
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -142858608
$ws.Range("N61").ClearContents()

# Row 122
$ws.Range("H122").Value = 5333.077
$ws.Range("I122").Value = 4824.316
$ws.Range("K122").Value = 14472.948
$ws.Range("M122").Value = -12022.948

# Row 132
$ws.Range("H132").Value = 2780551
$ws.Range("I132").Value = 2780551
$ws.Range("K132").Value = 8341653
$ws.Range("M132").Value = -8339123

# Row 136
$ws.Range("H136").Value = 142858820
$ws.Range("I136").Value = 142858820
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 428576460
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -428573910
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 13
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

# Row 107
$ws.Range("H107").Value = 638147.44
$ws.Range("I107").Value = 694929.6
$ws.Range("K107").Value = 694929.6
$ws.Range("M107").Value = -693009.6

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 6247.231
$ws.Range("I3").Value = 4655.909
$ws.Range("K3").Value = 13967.727
$ws.Range("M3").Value = -13855.727

# Row 33
$ws.Range("H33").Value = 319.82352
$ws.Range("I33").Value = 293.5
$ws.Range("K33").Value = 1761
$ws.Range("M33").Value = -1478

# Row 76
$ws.Range("H76").Value = 18686
$ws.Range("J76").Value = 18199.8
$ws.Range("L76").Value = 54599.39999999999
$ws.Range("N76").Value = -55365.39999999999

# Row 79
$ws.Range("H79").Value = 18686
$ws.Range("J79").Value = 18199.8
$ws.Range("L79").Value = 54599.39999999999
$ws.Range("N79").Value = -57251.39999999999

# Row 129
$ws.Range("H129").Value = 3366.7368
$ws.Range("J129").Value = 3748.2856
$ws.Range("L129").Value = 11244.8568
$ws.Range("N129").Value = -21244.8568

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4473.5
$ws.Range("I70").Value = 4473.5
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4473.5
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4203.5
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 4473.5
$ws.Range("I73").Value = 4473.5
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4473.5
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3537.5
$ws.Range("N73").ClearContents()

# Row 97
$ws.Range("H97").Value = 1362.4054
$ws.Range("I97").Value = 1285.2916
$ws.Range("K97").Value = 1285.2916
$ws.Range("M97").Value = -789.2916

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3368.3845
$ws.Range("I7").Value = 3149.75
$ws.Range("J7").Value = 3718.2
$ws.Range("K7").Value = 3149.75
$ws.Range("L7").Value = 3718.2
$ws.Range("M7").Value = -3037.75
$ws.Range("N7").Value = -3942.2

# Row 14
$ws.Range("H14").Value = 5004
$ws.Range("I14").Value = 5004
$ws.Range("K14").Value = 5004
$ws.Range("M14").Value = -4832

# Row 16
$ws.Range("H16").Value = 1756.7391
$ws.Range("I16").Value = 1190.5333
$ws.Range("J16").Value = 2818.375
$ws.Range("K16").Value = 1190.5333
$ws.Range("L16").Value = 2818.375
$ws.Range("M16").Value = -1020.5333
$ws.Range("N16").Value = -3158.375

# Row 40
$ws.Range("H40").Value = 2941.2354
$ws.Range("I40").Value = 2866.8
$ws.Range("K40").Value = 2866.8
$ws.Range("M40").Value = -2730.8

# Row 46
$ws.Range("H46").Value = 1617.9
$ws.Range("I46").Value = 1715.1765
$ws.Range("J46").Value = 1066.6666
$ws.Range("K46").Value = 1715.1765
$ws.Range("L46").Value = 1066.6666
$ws.Range("M46").Value = -1527.1765
$ws.Range("N46").Value = -1442.6666

# Row 68
$ws.Range("H68").Value = 6584946
$ws.Range("J68").Value = 7998.5
$ws.Range("L68").Value = 7998.5
$ws.Range("N68").Value = -9496.5

# Row 71
$ws.Range("H71").Value = 6584946
$ws.Range("J71").Value = 7998.5
$ws.Range("L71").Value = 39992.5
$ws.Range("N71").Value = -47480.5

# Row 122
$ws.Range("H122").Value = 5569.778
$ws.Range("I122").Value = 5000.5
$ws.Range("K122").Value = 15001.5
$ws.Range("M122").Value = -12551.5

# Row 126
$ws.Range("H126").Value = 3368.3845
$ws.Range("I126").Value = 3149.75
$ws.Range("J126").Value = 3718.2
$ws.Range("K126").Value = 9449.25
$ws.Range("L126").Value = 11154.6
$ws.Range("M126").Value = -6979.25
$ws.Range("N126").Value = -16094.6

# Row 136
$ws.Range("H136").Value = 2188.875
$ws.Range("I136").Value = 1760.75
$ws.Range("J136").Value = 2402.9375
$ws.Range("K136").Value = 5282.25
$ws.Range("L136").Value = 7208.8125
$ws.Range("M136").Value = -2732.25
$ws.Range("N136").Value = -12308.8125

$ws = $wb.Worksheets.Item("WVR")
# Row 29
$ws.Range("H29").Value = 9999.25
$ws.Range("J29").Value = 9999.25
$ws.Range("L29").Value = 9999.25
$ws.Range("N29").Value = -10579.25

# Row 80
$ws.Range("H80").Value = 79998
$ws.Range("J80").Value = 79998
$ws.Range("L80").Value = 79998
$ws.Range("N80").Value = -81994

# Row 83
$ws.Range("H83").Value = 79998
$ws.Range("J83").Value = 79998
$ws.Range("L83").Value = 239994
$ws.Range("N83").Value = -249978

# Row 103
$ws.Range("H103").Value = 28601
$ws.Range("J103").Value = 28601
$ws.Range("L103").Value = 28601
$ws.Range("N103").Value = -30945
